$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for price cells that would otherwise be parsed as numbers
$textCells = @("D4", "D5", "D7", "D8", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D18", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = '26.365.19'
$ws.Range("E2").Value = '  +3.50%  '
$ws.Range("D3").Value = '1.723.84'
$ws.Range("D4").Value = '0.9987'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '241.66'
$ws.Range("E5").Value = '  +1.70%  '
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("D7").Value = '0.4740'
$ws.Range("E7").Value = '  -1.18%  '
$ws.Range("D8").Value = '0.2638'
$ws.Range("E8").Value = '  +0.21%  '
$ws.Range("D9").Value = '0.06202'
$ws.Range("E9").Value = '  +0.32%  '
$ws.Range("D10").Value = '1.719.57'
$ws.Range("D11").Value = '0.07067'
$ws.Range("E11").Value = '  +0.02%  '
$ws.Range("D12").Value = '15.50'
$ws.Range("E12").Value = '  +4.40%  '
$ws.Range("D13").Value = '0.5976'
$ws.Range("E13").Value = '  +1.44%  '
$ws.Range("D14").Value = '4.426'
$ws.Range("E14").Value = '  +1.25%  '
$ws.Range("D15").Value = '76.39'
$ws.Range("E15").Value = '  +1.88%  '
$ws.Range("D16").Value = '0.9993'
$ws.Range("E16").Value = '  -0.10%  '
$ws.Range("D17").Value = '26.378.13'
$ws.Range("E17").Value = '  +3.54%  '
$ws.Range("D18").Value = '0.9992'
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("E19").Value = '  +1.07%  '
$ws.Range("D20").Value = '11.56'
$ws.Range("E20").Value = '  +0.97%  '
$ws.Range("D21").Value = '1.939.39'
$ws.Range("D22").Value = '4.526'
$ws.Range("E22").Value = '  +1.76%  '
$ws.Range("D23").Value = '8.754'
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").Value = '5.258'
$ws.Range("E24").Value = '  -0.54%  '
$ws.Range("D25").Value = '135.19'
$ws.Range("E25").Value = '  -1.24%  '
$ws.Range("D26").Value = '15.25'
$ws.Range("E26").Value = '  +1.46%  '
$ws.Range("D27").Value = '1.771'
$ws.Range("E27").Value = '  +2.85%  '
$ws.Range("E28").Value = '  +0.73%  '
$ws.Range("D29").Value = '107.01'
$ws.Range("E29").Value = '  +1.96%  '
$ws.Range("D30").Value = '3.964'
$ws.Range("E30").Value = '  +0.35%  '
$ws.Range("D31").Value = '3.689'
$ws.Range("E31").Value = '  +1.17%  '
$ws.Range("D32").Value = '0.07806'
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("D33").Value = '0.04499'
$ws.Range("E33").Value = '  +6.81%  '
$ws.Range("B34").Value = 'Frax'
$ws.Range("C34").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D34").Value = '0.9989'
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '2.613'
$ws.Range("E35").Value = '  +0.48%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '0.9830'
$ws.Range("E36").Value = '  +3.64%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '0.6234'
$ws.Range("E37").Value = '  +2.25%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '0.9411'
$ws.Range("E38").Value = '  +9.69%  '
$ws.Range("B39").Value = 'Quant'
$ws.Range("C39").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D39").Value = '113.75'
$ws.Range("E39").Value = '  +17.70%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '2.458'
$ws.Range("E40").Value = '  -5.32%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = '1.943'
$ws.Range("E41").Value = '  +4.96%  '
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").Value = '1.001'
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '5.671'
$ws.Range("E43").Value = '  +17.23%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '0.01489'
$ws.Range("E44").Value = '  +1.13%  '
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '0.3836'
$ws.Range("E45").Value = '  +1.85%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = '0.1189'
$ws.Range("E46").Value = '  +6.35%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").Value = '6.363'
$ws.Range("E47").Value = '  +2.57%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.05274'
$ws.Range("E48").Value = '  +0.38%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '7.835'
$ws.Range("E49").Value = '  +6.82%  '
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '30.43'
$ws.Range("E50").Value = '  +2.11%  '
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").Value = '0.3387'
$ws.Range("E51").Value = '  +1.68%  '
